$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated Cypher text for B2 (CasesTab row / query column) ---
# Adds an OPTIONAL MATCH for cohort + a trailing `Cohort` column to the RETURN clause.
$b2 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Oligodendroglioma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@
$b2 = $b2.TrimEnd("`r","`n")

# --- Updated Cypher text for B4 (FilesTab row / query column) ---
# Drops the trailing `Study Code` return column (and its preceding comma).
$b4 = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Oligodendroglioma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis 
'@
$b4 = $b4.TrimEnd("`r","`n")

# --- New StatQuery text for C2/C3/C4 ---
# Replaces the old "number_of_files / number_of_sample / ..." single-count query
# with a broader Programs/Studies/Cases/Samples/Case Files/Study Files rollup.
$cNew = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
  WHERE samp.specific_sample_pathology IN ["Oligodendroglioma"]  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$cNew = $cNew.TrimEnd("`r","`n")

# Apply the query-text edits.
$ws.Range("B2").Value = $b2
$ws.Range("B4").Value = $b4

# The StatQuery column (C) is identical across all three rows - point every row at
# the new rollup query (this also retires the old single-count query entirely).
$ws.Range("C2").Value = $cNew
$ws.Range("C3").Value = $cNew
$ws.Range("C4").Value = $cNew

# --- Row heights change (content grew/shrank after the query edits) ---
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210
